$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: D3 "FSK控制的慢1/4相位,TIM1_CH2" -> "FSK控制,TIM1_CH2"
$ws.Range("D3").Value = "FSK控制,TIM1_CH2"

# Row 4: F4 gets a new label "送FSK引脚" (D4 text itself stays "FSK控制,TIM1_CH3")
$ws.Range("F4").Value = "送FSK引脚"

# New rows 9-10, formatted like the existing data rows (border + left/vcenter align)
$ws.Range("A2:G2").Copy()
$ws.Range("A9:G10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Row 9: PA13 / DIO / SWDIO
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "PA13"
$ws.Range("C9").Value = "DIO"
$ws.Range("D9").Value = "SWDIO"

# Row 10: PA14 / DI / SWCLK
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "PA14"
$ws.Range("C10").Value = "DI"
$ws.Range("D10").Value = "SWCLK"

# Match the recorded UI selection state from the edit
$ws.Range("D13").Select() | Out-Null
